$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-04-18"

# Update the column header label for the April row
$ws.Range("A5").Value = "April (through 04-18)"

# Update April row (row 5) values
$ws.Range("C5").Value = 17
$ws.Range("D5").Value = 38
$ws.Range("E5").Value = 31
$ws.Range("F5").Value = 29
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 66
$ws.Range("I5").Value = 80

# Update Total row (row 6) values
$ws.Range("C6").Value = 145
$ws.Range("D6").Value = 227
$ws.Range("E6").Value = 228
$ws.Range("F6").Value = 139
$ws.Range("G6").Value = 238
$ws.Range("H6").Value = 489
$ws.Range("I6").Value = 515
